$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 64) of data to the bottom of the table in columns A:B,
# continuing the existing "index / isbn" list.
$ws.Cells.Item(64, 1).Value = 64
$ws.Cells.Item(64, 2).Value = 9781492087830

# Match the formatting of the row above (row 63) so the new row picks up the
# same cell styles (bordered/centered index column, numeric-formatted isbn
# column) instead of the default style.
$ws.Range("A63:B63").Copy()
$ws.Range("A64:B64").PasteSpecial(-4122)
$excel.CutCopyMode = $false
